# Ntrk3-Ptprf.xlsx update: new TPM-derived NATMI output.
# The original rows 2-4 (Sending cluster = "ECs") are no longer part of the
# output — the pipeline was re-run and "ECs" is not a sending cluster in the
# new result, only a target cluster. The remaining rows (old rows 5-10, for
# Sending cluster = "FAPs"/"MuSCs") move up to become rows 2-7, and all of
# their computed columns (E:T) are refreshed with newly recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows whose Sending cluster is "ECs" (old rows 2-4).
$ws.Rows("2:4").Delete()

# --- Row 2 (FAPs / Ntrk3 / Ptprf / ECs) ---
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ntrk3"
$ws.Range("C2").Value = "Ptprf"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.468673666666667
$ws.Range("H2").Value = 7.406021
$ws.Range("I2").Value = 0.635345274347677
$ws.Range("J2").Value = 0.635345274347677
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.05031533333333333
$ws.Range("N2").Value = 0.150946
$ws.Range("O2").Value = 0.005485022167780355
$ws.Range("P2").Value = 0.005485022167780356
$ws.Range("Q2").Value = 0.1242121384295555
$ws.Range("R2").Value = 1.117909245866
$ws.Range("S2").Value = 0.003484882913991499
$ws.Range("T2").Value = 0.0034848829139915

# --- Row 3 (FAPs / Ntrk3 / Ptprf / FAPs) ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ntrk3"
$ws.Range("C3").Value = "Ptprf"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.468673666666667
$ws.Range("H3").Value = 7.406021
$ws.Range("I3").Value = 0.635345274347677
$ws.Range("J3").Value = 0.635345274347677
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.467027333333334
$ws.Range("N3").Value = 10.401082
$ws.Range("O3").Value = 0.377950825718477
$ws.Range("P3").Value = 0.377950825718477
$ws.Range("Q3").Value = 8.558959079413556
$ws.Range("R3").Value = 77.03063171472201
$ws.Range("S3").Value = 0.2401292710560368
$ws.Range("T3").Value = 0.2401292710560368

# --- Row 4 (FAPs / Ntrk3 / Ptprf / MuSCs) ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ntrk3"
$ws.Range("C4").Value = "Ptprf"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.468673666666667
$ws.Range("H4").Value = 7.406021
$ws.Range("I4").Value = 0.635345274347677
$ws.Range("J4").Value = 0.635345274347677
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.655880666666666
$ws.Range("N4").Value = 16.967642
$ws.Range("O4").Value = 0.6165641521137426
$ws.Range("P4").Value = 0.6165641521137426
$ws.Range("Q4").Value = 13.96252366360911
$ws.Range("R4").Value = 125.662712972482
$ws.Range("S4").Value = 0.3917311203776486
$ws.Range("T4").Value = 0.3917311203776486

# --- Row 5 (MuSCs / Ntrk3 / Ptprf / ECs) ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Ntrk3"
$ws.Range("C5").Value = "Ptprf"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.416888666666667
$ws.Range("H5").Value = 4.250666
$ws.Range("I5").Value = 0.364654725652323
$ws.Range("J5").Value = 0.364654725652323
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.05031533333333333
$ws.Range("N5").Value = 0.150946
$ws.Range("O5").Value = 0.005485022167780355
$ws.Range("P5").Value = 0.005485022167780356
$ws.Range("Q5").Value = 0.07129122555955555
$ws.Range("R5").Value = 0.6416210300359999
$ws.Range("S5").Value = 0.002000139253788855
$ws.Range("T5").Value = 0.002000139253788856

# --- Row 6 (MuSCs / Ntrk3 / Ptprf / FAPs) ---
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Ntrk3"
$ws.Range("C6").Value = "Ptprf"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.416888666666667
$ws.Range("H6").Value = 4.250666
$ws.Range("I6").Value = 0.364654725652323
$ws.Range("J6").Value = 0.364654725652323
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.467027333333334
$ws.Range("N6").Value = 10.401082
$ws.Range("O6").Value = 0.377950825718477
$ws.Range("P6").Value = 0.377950825718477
$ws.Range("Q6").Value = 4.912391735623556
$ws.Range("R6").Value = 44.211525620612
$ws.Range("S6").Value = 0.1378215546624402
$ws.Range("T6").Value = 0.1378215546624402

# --- Row 7 (MuSCs / Ntrk3 / Ptprf / MuSCs) ---
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Ntrk3"
$ws.Range("C7").Value = "Ptprf"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.416888666666667
$ws.Range("H7").Value = 4.250666
$ws.Range("I7").Value = 0.364654725652323
$ws.Range("J7").Value = 0.364654725652323
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.655880666666666
$ws.Range("N7").Value = 16.967642
$ws.Range("O7").Value = 0.6165641521137426
$ws.Range("P7").Value = 0.6165641521137426
$ws.Range("Q7").Value = 8.01375321661911
$ws.Range("R7").Value = 72.12377894957199
$ws.Range("S7").Value = 0.224833031736094
$ws.Range("T7").Value = 0.224833031736094
